$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph-mark rFonts cleanup: the empty paragraph right before the
#    "CloudinaryDotNet..." paragraph currently carries an explicit
#    w:eastAsia="Times New Roman" override on its paragraph-mark rPr; the
#    target drops that attribute (ascii/hAnsi/cs stay untouched).
# ---------------------------------------------------------------------------
$cloudPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "CloudinaryDotNet*") {
        $cloudPara = $p
        break
    }
}
if ($cloudPara -eq $null) {
    throw "Could not locate the CloudinaryDotNet paragraph"
}

$emptyBeforeCloud = $cloudPara.Previous().Range
$fixPPrXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$emptyBeforeCloud.InsertXML($fixPPrXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) Drop the old (stray) _GoBack bookmark from its current paragraph near
#    the end of the Angular links list - it is relocated (see step 3 below)
#    to the newly-added paragraph block after the "Identity Package" entry.
# ---------------------------------------------------------------------------
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

# ---------------------------------------------------------------------------
# 3) Insert the three new paragraphs documenting
#    Microsoft.AspNetCore.Identity.EntityFrameworkCore right after the
#    "CloudinaryDotNet..." paragraph (and before the trio of blank
#    paragraphs that follow it). The middle paragraph carries the relocated
#    _GoBack bookmark.
# ---------------------------------------------------------------------------
$cloudPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "CloudinaryDotNet*") {
        $cloudPara = $p
        break
    }
}
if ($cloudPara -eq $null) {
    throw "Could not re-locate the CloudinaryDotNet paragraph"
}

$insertPoint = $d.Range($cloudPara.Range.End, $cloudPara.Range.End)
$newParasXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>
<w:p>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
      <w:color w:val="D4D4D4"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
      <w:b/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w14:shadow w14:blurRad="38100" w14:dist="19050" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="tl">
        <w14:schemeClr w14:val="dk1">
          <w14:alpha w14:val="60000"/>
        </w14:schemeClr>
      </w14:shadow>
      <w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr">
        <w14:noFill/>
        <w14:prstDash w14:val="solid"/>
        <w14:round/>
      </w14:textOutline>
    </w:rPr>
    <w:t>Microsoft.AspNetCore.Identity.EntityFrameworkCore</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
      <w:b/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w14:shadow w14:blurRad="38100" w14:dist="19050" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="tl">
        <w14:schemeClr w14:val="dk1">
          <w14:alpha w14:val="60000"/>
        </w14:schemeClr>
      </w14:shadow>
      <w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr">
        <w14:noFill/>
        <w14:prstDash w14:val="solid"/>
        <w14:round/>
      </w14:textOutline>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
      <w:bCs/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w14:shadow w14:blurRad="38100" w14:dist="19050" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="tl">
        <w14:schemeClr w14:val="dk1">
          <w14:alpha w14:val="60000"/>
        </w14:schemeClr>
      </w14:shadow>
      <w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr">
        <w14:noFill/>
        <w14:prstDash w14:val="solid"/>
        <w14:round/>
      </w14:textOutline>
    </w:rPr>
    <w:t>(Identity Package</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
      <w:bCs/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w14:shadow w14:blurRad="38100" w14:dist="19050" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="tl">
        <w14:schemeClr w14:val="dk1">
          <w14:alpha w14:val="60000"/>
        </w14:schemeClr>
      </w14:shadow>
      <w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr">
        <w14:noFill/>
        <w14:prstDash w14:val="solid"/>
        <w14:round/>
      </w14:textOutline>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
      <w:bCs/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w14:shadow w14:blurRad="38100" w14:dist="19050" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="tl">
        <w14:schemeClr w14:val="dk1">
          <w14:alpha w14:val="60000"/>
        </w14:schemeClr>
      </w14:shadow>
      <w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr">
        <w14:noFill/>
        <w14:prstDash w14:val="solid"/>
        <w14:round/>
      </w14:textOutline>
    </w:rPr>
    <w:t>By</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
      <w:bCs/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w14:shadow w14:blurRad="38100" w14:dist="19050" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="tl">
        <w14:schemeClr w14:val="dk1">
          <w14:alpha w14:val="60000"/>
        </w14:schemeClr>
      </w14:shadow>
      <w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr">
        <w14:noFill/>
        <w14:prstDash w14:val="solid"/>
        <w14:round/>
      </w14:textOutline>
    </w:rPr>
    <w:t xml:space="preserve"> Microsoft</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
      <w:bCs/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w14:shadow w14:blurRad="38100" w14:dist="19050" w14:dir="2700000" w14:sx="100000" w14:sy="100000" w14:kx="0" w14:ky="0" w14:algn="tl">
        <w14:schemeClr w14:val="dk1">
          <w14:alpha w14:val="60000"/>
        </w14:schemeClr>
      </w14:shadow>
      <w14:textOutline w14:w="0" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr">
        <w14:noFill/>
        <w14:prstDash w14:val="solid"/>
        <w14:round/>
      </w14:textOutline>
    </w:rPr>
    <w:t>)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
      <w:color w:val="D4D4D4"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
      <w:color w:val="D4D4D4"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insertPoint.InsertXML($newParasXml) | Out-Null

# ---------------------------------------------------------------------------
# 4) New character style "title" (based on Default Paragraph Font), appended
#    at the end of the style sheet.
# ---------------------------------------------------------------------------
$titleStyle = $d.Styles.Add("title", 2)
$titleStyle.BaseStyle = $d.Styles("DefaultParagraphFont")

Write-Output "done"
